$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AM")

for ($row = 2; $row -le 11; $row++) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").ClearContents()
    }
}
